$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels (order matters so shared-string table indices line up
# with the target: 0=SMallGrade, 1=BigGrade, 2=BossGrade, 3=Name)
$ws.Range("B1").Value = "SMallGrade"
$ws.Range("C1").Value = "BigGrade"
$ws.Range("D1").Value = "BossGrade"
$ws.Range("A1").Value = "Name"

# A1 should pick up the same "center/center" style already used by B1:D1 (xf index 1)
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108

# Remove the stray empty, styled-only cells in columns E:G for rows 1-9
$ws.Range("E1:G9").Clear()

# Give E13 a new "left/top" alignment style
$ws.Range("E13").HorizontalAlignment = -4131
$ws.Range("E13").VerticalAlignment = -4160

# Update the active selection shown when the sheet is opened
$ws.Range("F8").Select()
